$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "1.15", "36.696.37") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.696.37"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.080.84"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "244.72"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "53.60"
$ws.Range("E8").Value = "  -6.54%  "
$ws.Range("D9").Value = "58.63"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "0.365"
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").Value = "0.0760"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "14.94"
$ws.Range("E13").Value = "  -6.31%  "
$ws.Range("D14").Value = "0.884"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "2.386.58"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "2.097.10"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "36.667.42"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "17.22"
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("D20").Value = "72.47"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "0.0₃0875"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "5.42"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "240.32"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").Value = "9.83"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "167.22"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "20.67"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "5.28"
$ws.Range("E31").Value = "  +9.75%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  +5.10%  "
$ws.Range("D33").Value = "4.65"
$ws.Range("E33").Value = "  +3.99%  "
$ws.Range("D34").Value = "0.0605"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +6.08%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").Value = "0.0821"
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.15"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0219"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").Value = "4.82"
$ws.Range("E42").Value = "  -7.21%  "
$ws.Range("D43").Value = "0.0951"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "96.03"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.85"
$ws.Range("E45").Value = "  -10.66%  "
$ws.Range("D46").Value = "15.93"
$ws.Range("E46").Value = "  -7.84%  "
$ws.Range("D47").Value = "1.365.16"
$ws.Range("E47").Value = "  +7.14%  "
$ws.Range("D48").Value = "7.31"
$ws.Range("E48").Value = "  +7.33%  "
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "2.266.90"
$ws.Range("E51").Value = "  +1.54%  "
